# Hortaliza / Betarraga weekly update:
# A new week's data (row 92/93) is inserted at the top of the date-ordered
# block (rows 92-203). Every existing week shifts down by one pair of rows
# (Primera/Segunda), and the oldest week (which was at rows 202/203) is
# pushed out into two brand-new rows (204/205) at the end of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 92
$endRow = 203
$count = $endRow - $startRow + 1

# --- Read current (pre-edit) D/J/K/L/M/P columns for rows 92..203 ---
$Dvals = @()
$Jvals = @()
$Kvals = @()
$Lvals = @()
$Mvals = @()
$Pvals = @()

for ($i = 0; $i -lt $count; $i++) {
    $r = $startRow + $i
    $Dvals += $ws.Cells.Item($r, 4).Value2
    $Jvals += $ws.Cells.Item($r, 10).Value2
    $Kvals += $ws.Cells.Item($r, 11).Value2
    $Lvals += $ws.Cells.Item($r, 12).Value2
    $Mvals += $ws.Cells.Item($r, 13).Value2
    $Pvals += $ws.Cells.Item($r, 16).Value2
}

# --- Shift every pair (rows 94..203) down from two rows above it ---
for ($i = $count - 1; $i -ge 2; $i--) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 4).Value2 = $Dvals[$i - 2]
    $ws.Cells.Item($r, 10).Value2 = $Jvals[$i - 2]
    $ws.Cells.Item($r, 11).Value2 = $Kvals[$i - 2]
    $ws.Cells.Item($r, 12).Value2 = $Lvals[$i - 2]
    $ws.Cells.Item($r, 13).Value2 = $Mvals[$i - 2]
    $ws.Cells.Item($r, 16).Value2 = $Pvals[$i - 2]
}

# --- New week at the top: rows 92 (Primera) / 93 (Segunda) ---
# Dates move to 44539; volumes (J) change; min/max/avg price (K/L/M) and
# $/Kg (P) stay the same as they were before the shift.
$ws.Cells.Item(92, 4).Value2 = 44539
$ws.Cells.Item(92, 10).Value2 = 3000
$ws.Cells.Item(92, 11).Value2 = $Kvals[0]
$ws.Cells.Item(92, 12).Value2 = $Lvals[0]
$ws.Cells.Item(92, 13).Value2 = $Mvals[0]
$ws.Cells.Item(92, 16).Value2 = $Pvals[0]

$ws.Cells.Item(93, 4).Value2 = 44539
$ws.Cells.Item(93, 10).Value2 = 1400
$ws.Cells.Item(93, 11).Value2 = $Kvals[1]
$ws.Cells.Item(93, 12).Value2 = $Lvals[1]
$ws.Cells.Item(93, 13).Value2 = $Mvals[1]
$ws.Cells.Item(93, 16).Value2 = $Pvals[1]

# Make sure the new date cells keep the same date/time display format as
# the rest of column D.
$dateFmt = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(92, 4).NumberFormat = $dateFmt
$ws.Cells.Item(93, 4).NumberFormat = $dateFmt

# --- Oldest week (originally rows 202/203) is appended as new rows 204/205 ---
$newRows = @(204, 205)
$srcIdx = @($count - 2, $count - 1)   # indices for old rows 202, 203

for ($k = 0; $k -lt 2; $k++) {
    $r = $newRows[$k]
    $si = $srcIdx[$k]

    $ws.Cells.Item($r, 1).Value2 = 8
    $ws.Cells.Item($r, 2).Value2 = "Terminal La Palmera de La Serena"
    $ws.Cells.Item($r, 3).Value2 = "Coquimbo"
    $ws.Cells.Item($r, 4).Value2 = $Dvals[$si]
    $ws.Cells.Item($r, 4).NumberFormat = $dateFmt
    $ws.Cells.Item($r, 5).Value2 = 4
    $ws.Cells.Item($r, 6).Value2 = 100114014
    $ws.Cells.Item($r, 7).Value2 = "Betarraga"
    $ws.Cells.Item($r, 8).Value2 = "Sin especificar"
    if ($k -eq 0) {
        $ws.Cells.Item($r, 9).Value2 = "Primera"
    } else {
        $ws.Cells.Item($r, 9).Value2 = "Segunda"
    }
    $ws.Cells.Item($r, 10).Value2 = $Jvals[$si]
    $ws.Cells.Item($r, 11).Value2 = $Kvals[$si]
    $ws.Cells.Item($r, 12).Value2 = $Lvals[$si]
    $ws.Cells.Item($r, 13).Value2 = $Mvals[$si]
    $ws.Cells.Item($r, 14).Value2 = "`$/paquete 3 unidades"
    $ws.Cells.Item($r, 15).Value2 = "Provincia del Elqu" + [char]0x00ED
    $ws.Cells.Item($r, 16).Value2 = $Pvals[$si]
    $ws.Cells.Item($r, 17).Value2 = 3
    $ws.Cells.Item($r, 18).Value2 = "Hortaliza"
}
